$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Cluster'
$ws.Range("B1").Value = 'Active Cases'

$names = @(
  '3364 Assisi Centre Aged Care Rosanna',
  '3376 Royal Freemasons Coppin Centre Melbourne',
  '3622 Olivet Care Aged Care Services Ringwood',
  '3825 TLC Forest Lodge Residential Aged Care Frankston North',
  '3961 Heritage Care Water Gardens Aged Care Facility Sydenham',
  '4167 Royal Freemasons Centennial Lodge Wantirna South',
  '4282 Villa Maria Catholic Homes (VMCH) Wantirna At-Home Aged Care Stud Road Wantirna',
  '45573 Narre Warren South P-12 College Narre Warren South',
  '50567 Alamanda K9 College Point Cook',
  '52912 Edgars Creek Primary School Wollert',
  'AG Industries Pty Ltd Factory Thomastown',
  'Adass Israel School Elsternwick',
  'Antonine College Cedar Campus Coburg',
  'Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh',
  'Baden Powell College Tarneit',
  'Collingwood College Abbotsford',
  'Covenant College Bell Post Hill',
  'Dandenong South Primary School Dandenong',
  'Devon Meadows Primary School Devon Meadows',
  'Flemington Racecourse Flemington',
  'Gilly''s Early Learning Centre Balaclava',
  'Guardian Childcare & Education Moorabbin',
  'Hazel Glen College Doreen',
  'Hazelwood North Primary School Hazelwood North',
  'Ilim College Dallas Main Campus Dallas Oct',
  'Ilim College Glenroy Campus Hadfield',
  'Islamic College of Melbourne Tarneit Oct Nov',
  'Lyndhurst Primary School Lyndhurst',
  'Master Poultry Group West Footscray',
  'Minaret College Officer Campus Officer',
  'Morwell Park Primary School Morwell',
  'Nido Early School Woodend',
  'Nio Early Learning Adventures Preston',
  'Northern Bay College Wexford Campus Corio',
  'Northern Health Northern Hospital Epping Emergency Department Tier 1B',
  'Northern Health The Northern Hospital Epping',
  'Oakleigh South Primary School Oakleigh South',
  'Pentland Primary School Darley',
  'Rutherglen Motor Inn and Walkabout Motel Rutherglen',
  'Sirius College Ibrahim Dellal Campus Sunshine',
  'Sirius College Shepparton Campus Shepparton',
  'Smartie Pants Early Learning and Development Diamond Creek',
  'Social Gathering Woodvale 30 Oct',
  'Society Restaurant Melbourne',
  'St Ambrose Parish Primary School Woodend',
  'St Clare''s Primary School Officer',
  'St Georges Road Primary School Shepparton',
  'St Joseph''s School Quarry Hill',
  'St Louis de Montfort''s School Aspendale',
  'St Paul''s Primary School Sunshine West',
  'St Vincents Hospital Melbourne Emergency Department Fitzroy',
  'Stevensville Primary School St Albans',
  'Stockdale Road Primary School Traralgon',
  'Supreme Caravans Manufacturing Campbellfield',
  'Templestowe Park Primary School Templestowe',
  'The Lake Primary School Cabarita',
  'The Royal Children''s Hospital Melbourne Emergency Department Parkville',
  'Top Yard Rooftop Melbourne',
  'Truganina P-9 College Truganina',
  'Tucker Road Bentleigh Primary School Bentleigh',
  'Warragul Regional College Warragul',
  'Werribee Mercy Hospital Emergency Department',
  'Western Health Sunshine Hospital Emergency Department St Albans',
  'Wodonga Primary School Wodonga',
  'Wyndham Christian College Wyndham Vale',
  'Yeshivah College St Kilda East'
)

$counts = @(
  29,
  23,
  12,
  15,
  20,
  23,
  11,
  16,
  13,
  10,
  14,
  11,
  12,
  32,
  11,
  10,
  28,
  13,
  12,
  12,
  11,
  12,
  14,
  29,
  10,
  10,
  45,
  14,
  12,
  25,
  60,
  12,
  12,
  16,
  15,
  17,
  11,
  12,
  22,
  13,
  22,
  20,
  10,
  26,
  12,
  11,
  15,
  32,
  13,
  14,
  14,
  11,
  33,
  48,
  32,
  24,
  15,
  14,
  10,
  11,
  19,
  18,
  15,
  12,
  11,
  24
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $counts[$i]
}

Write-Host "Done. Rows written:" $names.Length
